$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Fgf16"
$ws.Cells.Item(2,3).Value = "Fgfr2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.1628603333333333
$ws.Cells.Item(2,8).Value = 0.488581
$ws.Cells.Item(2,9).Value = 0.06904471801498467
$ws.Cells.Item(2,10).Value = 0.06904471801498467
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.09434
$ws.Cells.Item(2,14).Value = 0.28302
$ws.Cells.Item(2,15).Value = 0.05191071108246543
$ws.Cells.Item(2,16).Value = 0.05191071108246543
$ws.Cells.Item(2,17).Value = 0.01536424384666667
$ws.Cells.Item(2,18).Value = 0.13827819462
$ws.Cells.Item(2,19).Value = 0.003584160408646165
$ws.Cells.Item(2,20).Value = 0.003584160408646165

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Fgf16"
$ws.Cells.Item(3,3).Value = "Fgfr2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.1628603333333333
$ws.Cells.Item(3,8).Value = 0.488581
$ws.Cells.Item(3,9).Value = 0.06904471801498467
$ws.Cells.Item(3,10).Value = 0.06904471801498467
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.9431116666666667
$ws.Cells.Item(3,14).Value = 2.829335
$ws.Cells.Item(3,15).Value = 0.5189484550226392
$ws.Cells.Item(3,16).Value = 0.5189484550226391
$ws.Cells.Item(3,17).Value = 0.1535954804038889
$ws.Cells.Item(3,18).Value = 1.382359323635
$ws.Cells.Item(3,19).Value = 0.03583064974135008
$ws.Cells.Item(3,20).Value = 0.03583064974135007

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Fgf16"
$ws.Cells.Item(4,3).Value = "Fgfr2"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.1628603333333333
$ws.Cells.Item(4,8).Value = 0.488581
$ws.Cells.Item(4,9).Value = 0.06904471801498467
$ws.Cells.Item(4,10).Value = 0.06904471801498467
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.7798996666666667
$ws.Cells.Item(4,14).Value = 2.339699
$ws.Cells.Item(4,15).Value = 0.4291408338948954
$ws.Cells.Item(4,16).Value = 0.4291408338948954
$ws.Cells.Item(4,17).Value = 0.1270147196798889
$ws.Cells.Item(4,18).Value = 1.143132477119
$ws.Cells.Item(4,19).Value = 0.02962990786498843
$ws.Cells.Item(4,20).Value = 0.02962990786498843

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Fgf16"
$ws.Cells.Item(5,3).Value = "Fgfr2"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.769244333333333
$ws.Cells.Item(5,8).Value = 5.307733
$ws.Cells.Item(5,9).Value = 0.7500720009247772
$ws.Cells.Item(5,10).Value = 0.7500720009247773
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.09434
$ws.Cells.Item(5,14).Value = 0.28302
$ws.Cells.Item(5,15).Value = 0.05191071108246543
$ws.Cells.Item(5,16).Value = 0.05191071108246543
$ws.Cells.Item(5,17).Value = 0.1669105104066666
$ws.Cells.Item(5,18).Value = 1.50219459366
$ws.Cells.Item(5,19).Value = 0.03893677093105285
$ws.Cells.Item(5,20).Value = 0.03893677093105286

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Fgf16"
$ws.Cells.Item(6,3).Value = "Fgfr2"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.769244333333333
$ws.Cells.Item(6,8).Value = 5.307733
$ws.Cells.Item(6,9).Value = 0.7500720009247772
$ws.Cells.Item(6,10).Value = 0.7500720009247773
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.9431116666666667
$ws.Cells.Item(6,14).Value = 2.829335
$ws.Cells.Item(6,15).Value = 0.5189484550226392
$ws.Cells.Item(6,16).Value = 0.5189484550226391
$ws.Cells.Item(6,17).Value = 1.668594971950556
$ws.Cells.Item(6,18).Value = 15.017354747555
$ws.Cells.Item(6,19).Value = 0.3892487060356528
$ws.Cells.Item(6,20).Value = 0.3892487060356527

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Fgf16"
$ws.Cells.Item(7,3).Value = "Fgfr2"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1.769244333333333
$ws.Cells.Item(7,8).Value = 5.307733
$ws.Cells.Item(7,9).Value = 0.7500720009247772
$ws.Cells.Item(7,10).Value = 0.7500720009247773
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.7798996666666667
$ws.Cells.Item(7,14).Value = 2.339699
$ws.Cells.Item(7,15).Value = 0.4291408338948954
$ws.Cells.Item(7,16).Value = 0.4291408338948954
$ws.Cells.Item(7,17).Value = 1.379833065818555
$ws.Cells.Item(7,18).Value = 12.418497592367
$ws.Cells.Item(7,19).Value = 0.3218865239580717
$ws.Cells.Item(7,20).Value = 0.3218865239580717

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Fgf16"
$ws.Cells.Item(8,3).Value = "Fgfr2"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.4266613333333333
$ws.Cells.Item(8,8).Value = 1.279984
$ws.Cells.Item(8,9).Value = 0.180883281060238
$ws.Cells.Item(8,10).Value = 0.180883281060238
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.09434
$ws.Cells.Item(8,14).Value = 0.28302
$ws.Cells.Item(8,15).Value = 0.05191071108246543
$ws.Cells.Item(8,16).Value = 0.05191071108246543
$ws.Cells.Item(8,17).Value = 0.04025123018666667
$ws.Cells.Item(8,18).Value = 0.36226107168
$ws.Cells.Item(8,19).Value = 0.009389779742766405
$ws.Cells.Item(8,20).Value = 0.009389779742766407

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Fgf16"
$ws.Cells.Item(9,3).Value = "Fgfr2"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.4266613333333333
$ws.Cells.Item(9,8).Value = 1.279984
$ws.Cells.Item(9,9).Value = 0.180883281060238
$ws.Cells.Item(9,10).Value = 0.180883281060238
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.9431116666666667
$ws.Cells.Item(9,14).Value = 2.829335
$ws.Cells.Item(9,15).Value = 0.5189484550226392
$ws.Cells.Item(9,16).Value = 0.5189484550226391
$ws.Cells.Item(9,17).Value = 0.4023892811822222
$ws.Cells.Item(9,18).Value = 3.62150353064
$ws.Cells.Item(9,19).Value = 0.09386909924563631
$ws.Cells.Item(9,20).Value = 0.09386909924563631

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Fgf16"
$ws.Cells.Item(10,3).Value = "Fgfr2"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.4266613333333333
$ws.Cells.Item(10,8).Value = 1.279984
$ws.Cells.Item(10,9).Value = 0.180883281060238
$ws.Cells.Item(10,10).Value = 0.180883281060238
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.7798996666666667
$ws.Cells.Item(10,14).Value = 2.339699
$ws.Cells.Item(10,15).Value = 0.4291408338948954
$ws.Cells.Item(10,16).Value = 0.4291408338948954
$ws.Cells.Item(10,17).Value = 0.3327530316462222
$ws.Cells.Item(10,18).Value = 2.994777284816
$ws.Cells.Item(10,19).Value = 0.07762440207183527
$ws.Cells.Item(10,20).Value = 0.07762440207183528
